$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.438.44'
$ws.Range("E2").Value = '  -1.32%  '

$ws.Range("D3").Value = '2.579.82'
$ws.Range("E3").Value = '  -2.06%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.93'
$ws.Range("E5").Value = '  -2.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.16'
$ws.Range("E6").Value = '  -1.23%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  -1.30%  '

$ws.Range("D9").Value = '2.578.80'
$ws.Range("E9").Value = '  -2.08%  '

$ws.Range("E10").Value = '  -4.04%  '

$ws.Range("E11").Value = '  +0.12%  '

$ws.Range("E12").Value = '  -1.85%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.15'
$ws.Range("E13").Value = '  -1.74%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.62'
$ws.Range("E14").Value = '  -4.20%  '

$ws.Range("E15").Value = '  -2.10%  '

$ws.Range("E16").Value = '  -3.14%  '

$ws.Range("D17").Value = '66.277.99'
$ws.Range("E17").Value = '  -1.37%  '

$ws.Range("D18").Value = '2.594.40'
$ws.Range("E18").Value = '  -1.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.38'
$ws.Range("E19").Value = '  -6.51%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.69'
$ws.Range("E20").Value = '  -5.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '350.37'
$ws.Range("E21").Value = '  -2.60%  '

$ws.Range("E22").Value = '  -3.23%  '

$ws.Range("E23").Value = '  -2.32%  '

$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("E25").Value = '  -4.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '68.76'
$ws.Range("E26").Value = '  -2.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.91'
$ws.Range("E27").Value = '  -9.14%  '

$ws.Range("D28").Value = '2.713.42'
$ws.Range("E28").Value = '  -2.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.34%  '

$ws.Range("D30").Value = '0.0₃0982'
$ws.Range("E30").Value = '  -3.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '529.59'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.04'
$ws.Range("E32").Value = '  +1.32%  '

$ws.Range("E33").Value = '  -3.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.84'
$ws.Range("E34").Value = '  -3.64%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.131'
$ws.Range("E35").Value = '  -3.98%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("E37").Value = '  -3.77%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '156.41'
$ws.Range("E38").Value = '  -0.65%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.70'
$ws.Range("E39").Value = '  -2.58%  '

$ws.Range("E40").Value = '  -2.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.28'
$ws.Range("E41").Value = '  +1.89%  '

$ws.Range("E42").Value = '  -1.87%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.09'
$ws.Range("E43").Value = '  -2.14%  '

$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("E45").Value = '  -2.49%  '

$ws.Range("D46").Value = '0.0₆0284'
$ws.Range("E46").Value = '  -5.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '148.44'
$ws.Range("E47").Value = '  -2.60%  '

$ws.Range("E48").Value = '  -4.11%  '

$ws.Range("E49").Value = '  -3.27%  '

$ws.Range("E50").Value = '  -2.05%  '

$ws.Range("E51").Value = '  -1.62%  '
